$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt -
# Zanahoria". It belongs chronologically right after the current row 573, so
# insert a fresh row at 574 (pushing the existing 574..682 block down to
# 575..683, preserving all of their data/formatting) and populate the new
# row with the new record's values.
$ws.Rows.Item(574).Insert()

$row = 574
$ws.Cells.Item($row, 1).Value2  = 4
$ws.Cells.Item($row, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value2  = "Los Lagos"
$ws.Cells.Item($row, 4).Value2  = 45258
$ws.Cells.Item($row, 5).Value2  = 10
$ws.Cells.Item($row, 6).Value2  = 100114013
$ws.Cells.Item($row, 7).Value2  = "Zanahoria"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 900
$ws.Cells.Item($row, 11).Value2 = 8500
$ws.Cells.Item($row, 12).Value2 = 9000
$ws.Cells.Item($row, 13).Value2 = 8750
$ws.Cells.Item($row, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item($row, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value2 = 438
$ws.Cells.Item($row, 17).Value2 = 20
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 1, 4).NumberFormat
